$d = $word.ActiveDocument

# Remove the sentence "The consolidation reported sufficient staff to meet this
# portion of the requirements." (which previously spanned several runs) while
# keeping the surrounding text - " requirements of Paragraph 45 of the HUD
# agreement." followed directly by " At the time of this ".
$old = "The consolidation reported sufficient staff to meet this portion of the requirements. "
$new = ""

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
